$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells D1, E1 - copy style from C1 (bold/border/centered) and set values
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:E1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4

# Updated C column values (rows 2-9)
$ws.Range("C2").Value = -5.022459607464143
$ws.Range("C3").Value = -1.15825194988682
$ws.Range("C4").Value = -0.07204406301364299
$ws.Range("C5").Value = -0.4144953840754857
$ws.Range("C6").Value = 0.01451842867919532
$ws.Range("C7").Value = 0.1067097157949464
$ws.Range("C8").Value = 0.1344667699115433
$ws.Range("C9").Value = 0.02767546902356237

# New D column values (rows 2-9)
$ws.Range("D2").Value = -4.691668997759643
$ws.Range("D3").Value = -1.157997093995135
$ws.Range("D4").Value = 0.02113037539333502
$ws.Range("D5").Value = -0.1530309736455852
$ws.Range("D6").Value = -0.0448712215926229
$ws.Range("D7").Value = 0.06782082594601035
$ws.Range("D8").Value = 0.0397297986559543
$ws.Range("D9").Value = 0.02414576792187023

# New E column values (rows 2-9)
$ws.Range("E2").Value = -4.32738296756667
$ws.Range("E3").Value = -1.136669764590692
$ws.Range("E4").Value = 0.09319397166410644
$ws.Range("E5").Value = 0.05433179954793064
$ws.Range("E6").Value = -0.09941504713009811
$ws.Range("E7").Value = 0.03020309278771556
$ws.Range("E8").Value = -0.05018014344453085
$ws.Range("E9").Value = 0.02109546198953799
